$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'25.705.93"
$ws.Range("E2").Value = "  -3.16%  "

# Row 3
$ws.Range("D3").Value = "'1.764.16"
$ws.Range("E3").Value = "  -4.18%  "

# Row 4
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  +0.39%  "

# Row 5
$ws.Range("D5").Value = "'236.83"
$ws.Range("E5").Value = "  -8.91%  "

# Row 6
$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "  +0.29%  "

# Row 7
$ws.Range("D7").Value = "'0.4935"
$ws.Range("E7").Value = "  -5.86%  "

# Row 8
$ws.Range("D8").Value = "'42.22"
$ws.Range("E8").Value = "  -5.77%  "

# Row 9
$ws.Range("D9").Value = "'0.2482"
$ws.Range("E9").Value = "  -21.84%  "

# Row 10
$ws.Range("D10").Value = "'0.06066"
$ws.Range("E10").Value = "  -10.72%  "

# Row 11
$ws.Range("D11").Value = "'1.774.78"
$ws.Range("E11").Value = "  -3.46%  "

# Row 12
$ws.Range("D12").Value = "'0.06647"
$ws.Range("E12").Value = "  -14.46%  "

# Row 13
$ws.Range("D13").Value = "'14.46"
$ws.Range("E13").Value = "  -22.83%  "

# Row 14
$ws.Range("D14").Value = "'0.6072"
$ws.Range("E14").Value = "  -22.42%  "

# Row 15
$ws.Range("D15").Value = "'77.99"
$ws.Range("E15").Value = "  -11.31%  "

# Row 16
$ws.Range("D16").Value = "'4.340"
$ws.Range("E16").Value = "  -13.52%  "

# Row 17
$ws.Range("D17").Value = "'1.003"
$ws.Range("E17").Value = "  +0.27%  "

# Row 18
$ws.Range("D18").Value = "'1.004"
$ws.Range("E18").Value = "  +0.33%  "

# Row 19
$ws.Range("D19").Value = "'25.732.49"
$ws.Range("E19").Value = "  -3.13%  "

# Row 20
$ws.Range("D20").Value = "'11.09"
$ws.Range("E20").Value = "  -19.98%  "

# Row 21
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "'0.000006288"
$ws.Range("E21").Value = "  -20.93%  "

# Row 22
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "'2.003.69"
$ws.Range("E22").Value = "  -3.33%  "

# Row 23
$ws.Range("D23").Value = "'3.841"
$ws.Range("E23").Value = "  -16.72%  "

# Row 24
$ws.Range("D24").Value = "'5.139"
$ws.Range("E24").Value = "  -14.07%  "

# Row 25
$ws.Range("D25").Value = "'8.016"
$ws.Range("E25").Value = "  -14.09%  "

# Row 26
$ws.Range("D26").Value = "'132.47"
$ws.Range("E26").Value = "  -7.13%  "

# Row 27
$ws.Range("D27").Value = "'1.859"
$ws.Range("E27").Value = "  -15.91%  "

# Row 28
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "'1.402"
$ws.Range("E28").Value = "  -16.37%  "

# Row 29
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'14.33"
$ws.Range("E29").Value = "  -15.21%  "

# Row 30
$ws.Range("D30").Value = "'99.49"
$ws.Range("E30").Value = "  -11.12%  "

# Row 31
$ws.Range("D31").Value = "'0.08200"
$ws.Range("E31").Value = "  -5.99%  "

# Row 32
$ws.Range("D32").Value = "'3.598"
$ws.Range("E32").Value = "  -13.72%  "

# Row 33
$ws.Range("D33").Value = "'1.003"
$ws.Range("E33").Value = "  +0.39%  "

# Row 34
$ws.Range("D34").Value = "'3.163"
$ws.Range("E34").Value = "  -22.41%  "

# Row 35
$ws.Range("D35").Value = "'0.04284"
$ws.Range("E35").Value = "  -12.40%  "

# Row 36
$ws.Range("D36").Value = "'2.612"
$ws.Range("E36").Value = "  -8.63%  "

# Row 37
$ws.Range("E37").Value = "  -10.07%  "

# Row 38
$ws.Range("D38").Value = "'0.6112"
$ws.Range("E38").Value = "  -15.47%  "

# Row 39
$ws.Range("D39").Value = "'2.708"
$ws.Range("E39").Value = "  -12.44%  "

# Row 40
$ws.Range("D40").Value = "'2.105"
$ws.Range("E40").Value = "  -5.57%  "

# Row 41
$ws.Range("D41").Value = "'1.003"
$ws.Range("E41").Value = "  +0.22%  "

# Row 42
$ws.Range("D42").Value = "'101.43"
$ws.Range("E42").Value = "  -8.00%  "

# Row 43
$ws.Range("D43").Value = "'0.01448"
$ws.Range("E43").Value = "  -17.10%  "

# Row 44
$ws.Range("D44").Value = "'0.7762"
$ws.Range("E44").Value = "  -13.40%  "

# Row 45
$ws.Range("D45").Value = "'0.3826"
$ws.Range("E45").Value = "  -20.50%  "

# Row 46
$ws.Range("D46").Value = "'5.168"
$ws.Range("E46").Value = "  -12.61%  "

# Row 47
$ws.Range("D47").Value = "'6.111"
$ws.Range("E47").Value = "  -20.04%  "

# Row 48
$ws.Range("D48").Value = "'0.05174"
$ws.Range("E48").Value = "  -11.33%  "

# Row 49
$ws.Range("D49").Value = "'52.28"
$ws.Range("E49").Value = "  -12.36%  "

# Row 50
$ws.Range("D50").Value = "'1.003"
$ws.Range("E50").Value = "  -0.10%  "

# Row 51
$ws.Range("D51").Value = "'1.001"
$ws.Range("E51").Value = "  +0.07%  "

